$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.112.07'
$ws.Range('E2').Value = '  -2.75%  '
$ws.Range('D3').Value = '2.995.43'
$ws.Range('E3').Value = '  -2.75%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.16'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.25'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -6.49%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D9').Value = '2.991.78'
$ws.Range('E9').Value = '  -2.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.147'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -6.10%  '
$ws.Range('E11').Value = '  -2.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.453'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('E13').Value = '  -4.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.27'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -7.13%  '
$ws.Range('E15').Value = '  +1.67%  '
$ws.Range('D16').Value = '3.484.08'
$ws.Range('E16').Value = '  -2.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.06'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.98%  '
$ws.Range('D18').Value = '62.004.03'
$ws.Range('E18').Value = '  -2.71%  '
$ws.Range('D19').Value = '2.992.53'
$ws.Range('E19').Value = '  -2.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '459.66'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -5.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.91'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.59%  '
$ws.Range('E22').Value = '  -4.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.42'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.15'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.21'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -10.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.22'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -5.67%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.96'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -7.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  -3.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.98'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -6.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.08'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -7.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '28.13'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.107'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -4.29%  '
$ws.Range('B35').Value = 'Mantle'
$ws.Range('C35').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.02'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -4.56%  '
$ws.Range('B36').Value = 'PEPE'
$ws.Range('C36').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D36').Value = '0.0₃0792'
$ws.Range('E36').Value = '  -4.05%  '
$ws.Range('E37').Value = '  -5.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.09'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -6.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '50.22'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.03%  '
$ws.Range('E40').Value = '  -1.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.86'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -12.04%  '
$ws.Range('E42').Value = '  +1.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '392.30'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -10.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0355'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.12%  '
$ws.Range('E45').Value = '  -6.99%  '
$ws.Range('D46').Value = '2.720.76'
$ws.Range('E46').Value = '  -4.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '36.66'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -7.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '127.93'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.108'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.49%  '
$ws.Range('E51').Value = '  -3.51%  '
